# Apply the diff: move stack-trace line numbers / object hash from the
# "before" version to the "after" version (library upgrade 3.0.0 -> 3.1.0),
# and insert a new "RunBefores" frame in the JUnit stack trace.

$d = $word.ActiveDocument
$content = $d.Content

function Replace-All($find, $replace) {
    $r = $d.Content
    $ok = $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARNING: replacement not found for: $find"
    }
}

# 1) Object hash code changed.
Replace-All "MTableImpl@2cccf134" "MTableImpl@5fb07347"

# 2) Line number updates (each of these strings is unique in the document,
#    except M2DocEvaluator.java:1216 which appears 3 times and always
#    becomes M2DocEvaluator.java:1239).
Replace-All "PaginationServices.java:178)" "PaginationServices.java:253)"
Replace-All "M2DocEvaluator.java:559)" "M2DocEvaluator.java:586)"
Replace-All "M2DocEvaluator.java:1216)" "M2DocEvaluator.java:1239)"
Replace-All "M2DocEvaluator.java:1425)" "M2DocEvaluator.java:1464)"
Replace-All "M2DocEvaluator.java:287)" "M2DocEvaluator.java:296)"
Replace-All "M2DocEvaluator.java:276)" "M2DocEvaluator.java:281)"
Replace-All "M2DocUtils.java:694)" "M2DocUtils.java:805)"
Replace-All "AbstractTemplatesTestSuite.java:480)" "AbstractTemplatesTestSuite.java:511)"
Replace-All "AbstractTemplatesTestSuite.java:389)" "AbstractTemplatesTestSuite.java:420)"

# 3) Insert a new stack frame line
#       at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)
#    right before the (unique) occurrence of RunAfters.java:27 that
#    immediately follows "ParentRunner$2.evaluate(ParentRunner.java:268)".
$old = "ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)"
$new = "ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)"
Replace-All $old $new

Write-Output "done"
